$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so numeric-looking
# strings (e.g. "0.9986", "253.03") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.131.18'
$ws.Range('E2').Value = '  +4.31%  '
$ws.Range('D3').Value = '1.907.87'
$ws.Range('E3').Value = '  +5.38%  '
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '253.03'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').Value = '0.9990'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '0.5089'
$ws.Range('E7').Value = '  +2.78%  '
$ws.Range('D8').Value = '45.22'
$ws.Range('E8').Value = '  +4.51%  '
$ws.Range('D9').Value = '0.3039'
$ws.Range('E9').Value = '  +9.16%  '
$ws.Range('D10').Value = '0.06823'
$ws.Range('D11').Value = '1.906.80'
$ws.Range('E11').Value = '  +5.29%  '
$ws.Range('D12').Value = '17.34'
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('D13').Value = '0.07327'
$ws.Range('E13').Value = '  +3.37%  '
$ws.Range('D14').Value = '0.6924'
$ws.Range('E14').Value = '  +7.00%  '
$ws.Range('D15').Value = '86.97'
$ws.Range('E15').Value = '  +3.16%  '
$ws.Range('D16').Value = '4.920'
$ws.Range('E16').Value = '  +4.75%  '
$ws.Range('D17').Value = '30.129.59'
$ws.Range('E17').Value = '  +4.22%  '
$ws.Range('D18').Value = '0.000008257'
$ws.Range('E18').Value = '  +12.30%  '
$ws.Range('E20').Value = '  +6.46%  '
$ws.Range('D21').Value = '2.153.32'
$ws.Range('E21').Value = '  +5.22%  '
$ws.Range('D22').Value = '0.9979'
$ws.Range('D23').Value = '4.825'
$ws.Range('E23').Value = '  +5.16%  '
$ws.Range('D24').Value = '5.746'
$ws.Range('E24').Value = '  +7.19%  '
$ws.Range('D25').Value = '9.301'
$ws.Range('E25').Value = '  +4.93%  '
$ws.Range('D26').Value = '147.63'
$ws.Range('E26').Value = '  +3.74%  '
$ws.Range('D27').Value = '135.35'
$ws.Range('E27').Value = '  +4.44%  '
$ws.Range('D28').Value = '17.13'
$ws.Range('E28').Value = '  +4.47%  '
$ws.Range('D29').Value = '2.005'
$ws.Range('E29').Value = '  +5.54%  '
$ws.Range('E30').Value = '  -0.86%  '
$ws.Range('D31').Value = '4.288'
$ws.Range('E31').Value = '  +3.21%  '
$ws.Range('D32').Value = '0.08851'
$ws.Range('E32').Value = '  +5.89%  '
$ws.Range('D33').Value = '4.008'
$ws.Range('E33').Value = '  +4.76%  '
$ws.Range('D34').Value = '0.05059'
$ws.Range('E34').Value = '  +1.74%  '
$ws.Range('D35').Value = '1.144'
$ws.Range('E35').Value = '  +4.17%  '
$ws.Range('D36').Value = '0.7258'
$ws.Range('E36').Value = '  +7.76%  '
$ws.Range('D37').Value = '2.686'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').Value = '2.276'
$ws.Range('E39').Value = '  -2.57%  '
$ws.Range('D40').Value = '0.9646'
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('D41').Value = '0.01696'
$ws.Range('E41').Value = '  +6.38%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '0.4319'
$ws.Range('E43').Value = '  +5.27%  '
$ws.Range('D44').Value = '105.09'
$ws.Range('E44').Value = '  +5.19%  '
$ws.Range('D45').Value = '0.9990'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '7.652'
$ws.Range('E46').Value = '  +6.86%  '
$ws.Range('D47').Value = '0.1281'
$ws.Range('E47').Value = '  +4.94%  '
$ws.Range('E48').Value = '  +4.17%  '
$ws.Range('D49').Value = '33.22'
$ws.Range('E49').Value = '  +4.69%  '
$ws.Range('D50').Value = '8.460'
$ws.Range('E50').Value = '  +3.74%  '
$ws.Range('D51').Value = '0.3831'
$ws.Range('E51').Value = '  +5.21%  '
